# Update "想去人数" (interest count, column F) figures to the latest
# scraped values, as regenerated for the gh-pages output at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8844
$ws1.Range("F8").Value = 6201
$ws1.Range("F9").Value = 608
$ws1.Range("F11").Value = 309
$ws1.Range("F12").Value = 9153
$ws1.Range("F13").Value = 10576
$ws1.Range("F15").Value = 1075
$ws1.Range("F16").Value = 4815
$ws1.Range("F18").Value = 406
$ws1.Range("F26").Value = 1161
$ws1.Range("F28").Value = 1990
$ws1.Range("F30").Value = 576
$ws1.Range("F31").Value = 2558
$ws1.Range("F34").Value = 1643
$ws1.Range("F40").Value = 3222
$ws1.Range("F41").Value = 4205
$ws1.Range("F49").Value = 4174

# Sheet "全部类型" (All types) - same events, different row layout
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8844
$ws4.Range("F8").Value = 6201
$ws4.Range("F9").Value = 608
$ws4.Range("F10").Value = 9153
$ws4.Range("F11").Value = 9153
$ws4.Range("F12").Value = 10576
$ws4.Range("F15").Value = 1075
$ws4.Range("F16").Value = 4815
$ws4.Range("F18").Value = 406
$ws4.Range("F25").Value = 1161
$ws4.Range("F28").Value = 1990
$ws4.Range("F30").Value = 576
$ws4.Range("F31").Value = 2558
